$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44676
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 4000
$ws.Range("L2").Value = 4500
$ws.Range("M2").Value = 4250
$ws.Range("P2").Value = 71

# Row 3
$ws.Range("D3").Value = 44657
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 5500
$ws.Range("M3").Value = 5250
$ws.Range("P3").Value = 88

# Row 5
$ws.Range("D5").Value = 44281
$ws.Range("K5").Value = 5500
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = 5750
$ws.Range("P5").Value = 96

# Row 6
$ws.Range("D6").Value = 44935
$ws.Range("K6").Value = 6000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 6500
$ws.Range("P6").Value = 108

# Row 7
$ws.Range("D7").Value = 44785
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7500
$ws.Range("P7").Value = 125

# Row 8
$ws.Range("D8").Value = 44362
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 8500
$ws.Range("P8").Value = 142

# Row 9
$ws.Range("D9").Value = 44400
$ws.Range("J9").Value = 120
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 9500
$ws.Range("P9").Value = 158

# Row 10
$ws.Range("D10").Value = 44421
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 8500
$ws.Range("P10").Value = 142

# Row 11
$ws.Range("D11").Value = 44603
$ws.Range("J11").Value = 140
$ws.Range("K11").Value = 5500
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = 5750
$ws.Range("P11").Value = 96

# Row 12
$ws.Range("D12").Value = 44740
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 6000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 6500
$ws.Range("P12").Value = 108

# Row 13
$ws.Range("D13").Value = 44764
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 8000
$ws.Range("M13").Value = 7500
$ws.Range("P13").Value = 125

# Row 14
$ws.Range("D14").Value = 44669
$ws.Range("J14").Value = 130
$ws.Range("K14").Value = 4500
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = 4750
$ws.Range("P14").Value = 79

# Row 15
$ws.Range("D15").Value = 44589
$ws.Range("J15").Value = 110
$ws.Range("K15").Value = 5000
$ws.Range("M15").Value = 5500
$ws.Range("P15").Value = 92

# Row 16
$ws.Range("D16").Value = 44242
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 5500
$ws.Range("M16").Value = 5250
$ws.Range("P16").Value = 88

# Row 17
$ws.Range("D17").Value = 44494
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 5000
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = 5500
$ws.Range("P17").Value = 92

# Row 18
$ws.Range("D18").Value = 44760
$ws.Range("J18").Value = 130
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 7500
$ws.Range("M18").Value = 7250
$ws.Range("P18").Value = 121

# Row 19
$ws.Range("D19").Value = 44627
$ws.Range("K19").Value = 4000
$ws.Range("L19").Value = 4500
$ws.Range("M19").Value = 4250
$ws.Range("P19").Value = 71

# Row 21
$ws.Range("D21").Value = 44827
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 6000
$ws.Range("L21").Value = 7000
$ws.Range("M21").Value = 6500
$ws.Range("P21").Value = 108
